$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value2 = 17.70643966666667
$ws.Cells.Item(2, 8).Value2 = 53.119319
$ws.Cells.Item(2, 9).Value2 = 0.4380235920947999
$ws.Cells.Item(2, 10).Value2 = 0.4380235920947999
$ws.Cells.Item(2, 13).Value2 = 16.28844733333333
$ws.Cells.Item(2, 14).Value2 = 48.865342
$ws.Cells.Item(2, 15).Value2 = 0.2176904746803693
$ws.Cells.Item(2, 16).Value2 = 0.2176904746803693
$ws.Cells.Item(2, 17).Value2 = 288.4104099713443
$ws.Cells.Item(2, 18).Value2 = 2595.693689742098
$ws.Cells.Item(2, 19).Value2 = 0.09535356368431745
$ws.Cells.Item(2, 20).Value2 = 0.09535356368431745

$ws.Cells.Item(3, 7).Value2 = 17.70643966666667
$ws.Cells.Item(3, 8).Value2 = 53.119319
$ws.Cells.Item(3, 9).Value2 = 0.4380235920947999
$ws.Cells.Item(3, 10).Value2 = 0.4380235920947999
$ws.Cells.Item(3, 13).Value2 = 27.61090666666666
$ws.Cells.Item(3, 14).Value2 = 82.83272
$ws.Cells.Item(3, 15).Value2 = 0.3690119294748028
$ws.Cells.Item(3, 16).Value2 = 0.3690119294748029
$ws.Cells.Item(3, 17).Value2 = 488.8908530352978
$ws.Cells.Item(3, 18).Value2 = 4400.01767731768
$ws.Cells.Item(3, 19).Value2 = 0.1616359308743861
$ws.Cells.Item(3, 20).Value2 = 0.1616359308743861

$ws.Cells.Item(4, 7).Value2 = 17.70643966666667
$ws.Cells.Item(4, 8).Value2 = 53.119319
$ws.Cells.Item(4, 9).Value2 = 0.4380235920947999
$ws.Cells.Item(4, 10).Value2 = 0.4380235920947999
$ws.Cells.Item(4, 13).Value2 = 26.266325
$ws.Cells.Item(4, 14).Value2 = 78.798975
$ws.Cells.Item(4, 15).Value2 = 0.3510419771967738
$ws.Cells.Item(4, 16).Value2 = 0.3510419771967739
$ws.Cells.Item(4, 17).Value2 = 465.0830988775584
$ws.Cells.Item(4, 18).Value2 = 4185.747889898025
$ws.Cells.Item(4, 19).Value2 = 0.1537646678277917
$ws.Cells.Item(4, 20).Value2 = 0.1537646678277917

$ws.Cells.Item(5, 7).Value2 = 17.70643966666667
$ws.Cells.Item(5, 8).Value2 = 53.119319
$ws.Cells.Item(5, 9).Value2 = 0.4380235920947999
$ws.Cells.Item(5, 10).Value2 = 0.4380235920947999
$ws.Cells.Item(5, 13).Value2 = 4.658207333333333
$ws.Cells.Item(5, 14).Value2 = 13.974622
$ws.Cells.Item(5, 15).Value2 = 0.06225561864805391
$ws.Cells.Item(5, 16).Value2 = 0.06225561864805392
$ws.Cells.Item(5, 17).Value2 = 82.4802671024909
$ws.Cells.Item(5, 18).Value2 = 742.3224039224181
$ws.Cells.Item(5, 19).Value2 = 0.02726942970830458
$ws.Cells.Item(5, 20).Value2 = 0.02726942970830459

$ws.Cells.Item(6, 7).Value2 = 1.617245333333334
$ws.Cells.Item(6, 8).Value2 = 4.851736000000001
$ws.Cells.Item(6, 9).Value2 = 0.04000756919748267
$ws.Cells.Item(6, 10).Value2 = 0.04000756919748267
$ws.Cells.Item(6, 13).Value2 = 16.28844733333333
$ws.Cells.Item(6, 14).Value2 = 48.865342
$ws.Cells.Item(6, 15).Value2 = 0.2176904746803693
$ws.Cells.Item(6, 16).Value2 = 0.2176904746803693
$ws.Cells.Item(6, 17).Value2 = 26.34241543707912
$ws.Cells.Item(6, 18).Value2 = 237.081738933712
$ws.Cells.Item(6, 19).Value2 = 0.008709266729407725
$ws.Cells.Item(6, 20).Value2 = 0.008709266729407725

$ws.Cells.Item(7, 7).Value2 = 1.617245333333334
$ws.Cells.Item(7, 8).Value2 = 4.851736000000001
$ws.Cells.Item(7, 9).Value2 = 0.04000756919748267
$ws.Cells.Item(7, 10).Value2 = 0.04000756919748267
$ws.Cells.Item(7, 13).Value2 = 27.61090666666666
$ws.Cells.Item(7, 14).Value2 = 82.83272
$ws.Cells.Item(7, 15).Value2 = 0.3690119294748028
$ws.Cells.Item(7, 16).Value2 = 0.3690119294748029
$ws.Cells.Item(7, 19).Value2 = 0.01476327030315977
$ws.Cells.Item(7, 20).Value2 = 0.01476327030315977

$ws.Cells.Item(8, 7).Value2 = 1.617245333333334
$ws.Cells.Item(8, 8).Value2 = 4.851736000000001
$ws.Cells.Item(8, 9).Value2 = 0.04000756919748267
$ws.Cells.Item(8, 10).Value2 = 0.04000756919748267
$ws.Cells.Item(8, 13).Value2 = 26.266325
$ws.Cells.Item(8, 14).Value2 = 78.798975
$ws.Cells.Item(8, 15).Value2 = 0.3510419771967738
$ws.Cells.Item(8, 16).Value2 = 0.3510419771967739
$ws.Cells.Item(8, 17).Value2 = 42.47909153006668
$ws.Cells.Item(8, 18).Value2 = 382.3118237706
$ws.Cells.Item(8, 19).Value2 = 0.01404433619392106
$ws.Cells.Item(8, 20).Value2 = 0.01404433619392106

$ws.Cells.Item(9, 7).Value2 = 1.617245333333334
$ws.Cells.Item(9, 8).Value2 = 4.851736000000001
$ws.Cells.Item(9, 9).Value2 = 0.04000756919748267
$ws.Cells.Item(9, 10).Value2 = 0.04000756919748267
$ws.Cells.Item(9, 13).Value2 = 4.658207333333333
$ws.Cells.Item(9, 14).Value2 = 13.974622
$ws.Cells.Item(9, 15).Value2 = 0.06225561864805391
$ws.Cells.Item(9, 16).Value2 = 0.06225561864805392
$ws.Cells.Item(9, 17).Value2 = 7.533464071532446
$ws.Cells.Item(9, 18).Value2 = 67.801176643792
$ws.Cells.Item(9, 19).Value2 = 0.00249069597099411
$ws.Cells.Item(9, 20).Value2 = 0.00249069597099411

$ws.Cells.Item(10, 7).Value2 = 21.099799
$ws.Cells.Item(10, 8).Value2 = 63.299397
$ws.Cells.Item(10, 9).Value2 = 0.5219688387077175
$ws.Cells.Item(10, 10).Value2 = 0.5219688387077175
$ws.Cells.Item(10, 13).Value2 = 16.28844733333333
$ws.Cells.Item(10, 14).Value2 = 48.865342
$ws.Cells.Item(10, 15).Value2 = 0.2176904746803693
$ws.Cells.Item(10, 16).Value2 = 0.2176904746803693
$ws.Cells.Item(10, 17).Value2 = 343.6829647554193
$ws.Cells.Item(10, 18).Value2 = 3093.146682798774
$ws.Cells.Item(10, 19).Value2 = 0.1136276442666442
$ws.Cells.Item(10, 20).Value2 = 0.1136276442666442

$ws.Cells.Item(11, 7).Value2 = 21.099799
$ws.Cells.Item(11, 8).Value2 = 63.299397
$ws.Cells.Item(11, 9).Value2 = 0.5219688387077175
$ws.Cells.Item(11, 10).Value2 = 0.5219688387077175
$ws.Cells.Item(11, 13).Value2 = 27.61090666666666
$ws.Cells.Item(11, 14).Value2 = 82.83272
$ws.Cells.Item(11, 15).Value2 = 0.3690119294748028
$ws.Cells.Item(11, 16).Value2 = 0.3690119294748029
$ws.Cells.Item(11, 17).Value2 = 582.5845808744267
$ws.Cells.Item(11, 18).Value2 = 5243.26122786984
$ws.Cells.Item(11, 19).Value2 = 0.192612728297257
$ws.Cells.Item(11, 20).Value2 = 0.192612728297257

$ws.Cells.Item(12, 7).Value2 = 21.099799
$ws.Cells.Item(12, 8).Value2 = 63.299397
$ws.Cells.Item(12, 9).Value2 = 0.5219688387077175
$ws.Cells.Item(12, 10).Value2 = 0.5219688387077175
$ws.Cells.Item(12, 13).Value2 = 26.266325
$ws.Cells.Item(12, 14).Value2 = 78.798975
$ws.Cells.Item(12, 15).Value2 = 0.3510419771967738
$ws.Cells.Item(12, 16).Value2 = 0.3510419771967739
$ws.Cells.Item(12, 17).Value2 = 554.214177968675
$ws.Cells.Item(12, 18).Value2 = 4987.927601718075
$ws.Cells.Item(12, 19).Value2 = 0.1832329731750611
$ws.Cells.Item(12, 20).Value2 = 0.1832329731750611

$ws.Cells.Item(13, 7).Value2 = 21.099799
$ws.Cells.Item(13, 8).Value2 = 63.299397
$ws.Cells.Item(13, 9).Value2 = 0.5219688387077175
$ws.Cells.Item(13, 10).Value2 = 0.5219688387077175
$ws.Cells.Item(13, 13).Value2 = 4.658207333333333
$ws.Cells.Item(13, 14).Value2 = 13.974622
$ws.Cells.Item(13, 15).Value2 = 0.06225561864805391
$ws.Cells.Item(13, 16).Value2 = 0.06225561864805392
$ws.Cells.Item(13, 17).Value2 = 98.28723843365934
$ws.Cells.Item(13, 18).Value2 = 884.585145902934
$ws.Cells.Item(13, 19).Value2 = 0.03249549296875522
$ws.Cells.Item(13, 20).Value2 = 0.03249549296875523
